$wb = $excel.ActiveWorkbook

# The localization status changed from "Ready for handoff" to "In Translation"
# for the Overview sheet (zh-cn / de-de status columns) as well as the
# per-locale "Status" column on the "zh-cn" and "de-de" sheets.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# Shorter status text -> narrower auto-fitted "Status" columns. Re-applying
# the fitted width on every sheet that shows the status value.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
